$wb = $excel.ActiveWorkbook

# Metadata sheet
$ws = $wb.Worksheets.Item("Metadata")

# Update Publisher value (row 9, column B) - German -> English
$ws.Range("B9").Value = "Independent Trusted Third Party of the University Medicine Greifswald"

# Update Contact value (row 10, column B) - German -> English
$ws.Range("B10").Value = "Independent Trusted Third Party of the University Medicine Greifswald (https://www.ths-greifswald.de/)"

# Set Description value (row 12, column B) which was previously empty
$ws.Range("B12").Value = "consent states - minimal subset REFUSAL documents"
